$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.025.22'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.619.93'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").Value = '1.626.68'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = '27.003.00'
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.59%  '
$ws.Range("D17").Value = '0.0₃0741'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("E22").Value = '  -5.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '147.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.70%  '
$ws.Range("E27").Value = '  -0.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0504'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("E31").Value = '  -1.40%  '
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("D33").Value = '1.337.89'
$ws.Range("E33").Value = '  +6.00%  '
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.847'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.48%  '
$ws.Range("D43").Value = '1.756.78'
$ws.Range("E43").Value = '  -1.51%  '
$ws.Range("E44").Value = '  -2.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.851'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +28.10%  '
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0998'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0512'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("E51").Value = '  -1.40%  '
